$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 181
$ws.Range("I2").Value = 518
$ws.Range("J2").Value = 1906
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 539
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 345
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 21
$ws.Range("S2").Value = 198
$ws.Range("T2").Value = 340
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 3202
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2955
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 43
$ws.Range("AA2").Value = 15
